$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated measurements (fixed data + added extra frame cuts).
# Columns: C=Frame, D=I_max, E=I_min1, F=I_min2
$data = @{
    2  = @(3,   214.3, 42.4,              38.9)
    3  = @(7,   243.2, 42.3,              38.299999999999997)
    4  = @(11,  231.8, 43.7,              48.7)
    5  = @(19,  250.2, 42.6,              28.3)
    6  = @(24,  249.4, 50,                71.5)
    7  = @(46,  248.4, 43.2,              40.4)
    8  = @(62,  231.3, 59.8,              50.6)
    9  = @(95,  234.5, 56.5,              77.3)
    10 = @(122, 252.6, 49.5,              48.1)
    11 = @(176, 249.3, 35.799999999999997, 35)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}

# Move the selection like the author left it before saving.
$ws.Range("C12").Select()
